$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New raw-data rows appended after the existing data (rows 212-231),
# matching the T4 ("Dark", T3 timepoint label) PAM readings added upstream.
$rowsData = @(
    @{ Row=212; C=1; B="IPa06-C-T4"; F=36; H="Control"; I=30; L=0.583 }
    @{ Row=213; C=2; B="IPa07-C-T4"; F=16; H="Control"; I=30; L=0.613 }
    @{ Row=214; C=3; B="IPa08-C-T4"; F=38; H="Control"; I=30; L=0.623 }
    @{ Row=215; C=4; B="IPa09-C-T4"; F=33; H="Control"; I=30; L=0.603 }
    @{ Row=216; C=5; B="IPa10-C-T4"; F=6; H="Control"; I=30; L=0.549 }
    @{ Row=217; C=1; B="IPa06-L-T4"; F=36; H="Low"; I=34; L=0.613 }
    @{ Row=218; C=2; B="IPa07-L-T4"; F=16; H="Low"; I=34; L=0.596 }
    @{ Row=219; C=3; B="IPa08-L-T4"; F=38; H="Low"; I=34; L=0.62 }
    @{ Row=220; C=4; B="IPa09-L-T4"; F=33; H="Low"; I=34; L=0.604 }
    @{ Row=221; C=5; B="IPa10-L-T4"; F=6; H="Low"; I=34; L=0.591 }
    @{ Row=222; C=1; B="IPa06-M-T4"; F=36; H="Medium"; I=37; L=0.429 }
    @{ Row=223; C=2; B="IPa07-M-T4"; F=16; H="Medium"; I=37; L=0.533 }
    @{ Row=224; C=3; B="IPa08-M-T4"; F=38; H="Medium"; I=37; L=0.539 }
    @{ Row=225; C=4; B="IPa09-M-T4"; F=33; H="Medium"; I=37; L=0.349 }
    @{ Row=226; C=5; B="IPa10-M-T4"; F=6; H="Medium"; I=37; L=0.408 }
    @{ Row=227; C=1; B="IPa06-H-T4"; F=36; H="High"; I=39; L=0.051 }
    @{ Row=228; C=2; B="IPa07-H-T4"; F=16; H="High"; I=39; L=0.128 }
    @{ Row=229; C=3; B="IPa08-H-T4"; F=38; H="High"; I=39; L=0.069 }
    @{ Row=230; C=4; B="IPa09-H-T4"; F=33; H="High"; I=39; L=0.104 }
    @{ Row=231; C=5; B="IPa10-H-T4"; F=6; H="High"; I=39; L=0.027 }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row
    $ws.Cells.Item($r, 1).Value2 = "NA"             # TubeNumber
    $ws.Cells.Item($r, 2).Value2 = $rd.B             # SampleName
    $ws.Cells.Item($r, 3).Value2 = $rd.C             # Genotype
    $ws.Cells.Item($r, 4).Value2 = "Porites astreoides" # Species
    $ws.Cells.Item($r, 5).Value2 = "Inshore"         # Location
    $ws.Cells.Item($r, 6).Value2 = $rd.F             # TagNumber
    $ws.Cells.Item($r, 7).Value2 = "T3"              # Timepoint
    $ws.Cells.Item($r, 8).Value2 = $rd.H             # Treatment
    $ws.Cells.Item($r, 9).Value2 = $rd.I             # Temp_Setpoint
    $ws.Cells.Item($r, 10).Value2 = "Dark"           # PAM_type
    $ws.Cells.Item($r, 12).Value2 = $rd.L            # PAM
    $ws.Cells.Item($r, 13).Value2 = "NA"             # Time_RNAseq
    $ws.Cells.Item($r, 14).Value2 = 44355            # Date
}

# Carry over the same cell formatting used by the rest of the table
# (style used for columns A-J and M-N); column K (Time_PAM) and column O
# (Comments) are intentionally left blank/untouched for these rows, exactly
# as in the source diff.
$ws.Range("A211:J211").Copy()
$ws.Range("A212:J231").PasteSpecial(-4122)
$ws.Range("M211:N211").Copy()
$ws.Range("M212:N231").PasteSpecial(-4122)

# The PAM (L) column for this newly appended block uses a distinct style:
# a 12pt Calibri font in the theme text color, right aligned - not yet
# present in the workbook. Seed the format from the closest existing cell
# that already matches font/size/alignment (L22, which is RGB black rather
# than theme-colored), then switch the font color to the theme color so
# exactly one new font + one new cell style are produced.
$ws.Range("L22").Copy()
$ws.Range("L212").PasteSpecial(-4122)
$ws.Cells.Item(212, 12).Value2 = 0.583
$ws.Range("L212").Font.ThemeColor = 1

$ws.Range("L212").Copy()
$ws.Range("L213:L231").PasteSpecial(-4122)
foreach ($rd in $rowsData) {
    $ws.Cells.Item($rd.Row, 12).Value2 = $rd.L
}

$excel.CutCopyMode = $false

Write-Host "Appended 20 rows (212-231) of Inshore PAM raw data"
